# Update "想去人数" (number of people interested) values on both the
# "展览" (Exhibitions) and "全部类型" (All types) sheets, as produced by a
# fresh scrape/regeneration of the site's data (gh-pages output update).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8
$ws1.Range("F5").Value = 2552
$ws1.Range("F6").Value = 1841
$ws1.Range("F8").Value = 108
$ws1.Range("F9").Value = 896
$ws1.Range("F10").Value = 176

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8
$ws4.Range("F5").Value = 2552
$ws4.Range("F6").Value = 1841
$ws4.Range("F9").Value = 108
$ws4.Range("F10").Value = 896
$ws4.Range("F11").Value = 176
